$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAyoyYUNnHlwFT5mT6bpZe\nKbDtalgVOrp7YnJnRetXX8LmN/TpWQ54uEoK7qB6xFqFKrETTUDtoOoVbUer1U9/\nzeBFI+DG34LwL/b1LE8ht9FFy90MLLAO1vjnW2NCYZ9FT+7wpEB3Y25y694lcOAX\ndzPyeeEwAlAyLJtsvBcQQKfraop34Tmyf7VcySJRuxtmzQ1OWcb+HNqJxoXvceKZ\nIoHklgbUdEgVzFXEu1q7qr/P0wy29NlZGljtIPM4uhnd1T2O5O3e3w7Arem2Ie0Z\nbpVZcjvK+I8wAkANn1PF4Jd3OMI2Y1UQ+FFK6nPwtg/hapgtmDeIgJKmTc1RImF3\nMQIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D3").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAlUW3DMG89rUZHd/Be3ge\nsfmUreJBoUfFwX+fh0yO9rI0kohVvlkE1mRxK1yjzhWFn7iYJp/empNuhYwq1e6f\n/ecmvXCokDRbkePvj39fyYMv+4YpSNDTFemAS60+bL4/vCVzjo6UuCzyHQwbAodr\nPLgD9yG1qoiBPZdRMlzZ63WNZAwfaHwvB0f+1TPQ4V6GXJbpiUbsbPOhXF06IOCT\nr8xQqz/mMW8iCm227hKvmJwGXC0oCgPeDsNF3Ohc7p/lKwjxWpO363nyI7z4Z9QI\np7nZRRZcz6N9h8ulDPf9OgnZ7OCane/ddF+25hzO0R7Z6sWUAQFKLT4TtuvQLBLJ\nvwIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D4").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAz+jkAbxEeU0aaFzPUX9G\nKC8E/YmBDFAkLVroF9zGf1MbFBYe+U0kDQgtxHjOokAufGW6bWX7iDyU7pcUe+de\nRJe6PSs+WNcMt7bT+TnFrBOEusPSV3Z9a/ijX4j7HDCK235D2Xn1mo6dYpxRpHEg\nGKpilNDX8owuJZrGQla4X0a0LEvmHy82yQqg1OcMu9JDj5ZdULGC1yW3SQVdqMmW\nXpe4cfD/YyIXEw6otdaQanB9vIMwEIym08Bt2UtFhMy+hLtg+ui3VM2Foot62QPx\nRMJvi+KzNXzWhl16OvVsoXjLysH94Lg6eZD/OuCfo+C4gljrCe7k0o5V4297RwfY\nJwIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D5").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAhR2WP83Q6uDSSeaF8GvW\nnRzuzXce0B5tUmhlWpcGoQ0LFPi1kOnYKUQLWkPPhAOKqI0Q7sG5dGgGUvNVIEPf\nLIHFthY5EZoC1O7s9WX74gKWozo9gmaWuGaqS9Jk2+wO3KTPp1xdmVCsZYiWCb2d\nq3tCUqi5pbBKO3RwBonGoEL6vYbW8KE0TgokFYJDdhD69T2DyG9LT8kGhS10VTiI\nzXLXKtsoJRikCzfqSL8LXp8XulCPDyi0uGxCQNPUY170O17foW61/fyac3JlEAXf\nDDydJQZpJshky9hb/EfQxRx4BhDSbYVvHKloEf6NNTl6RKAhjs3+oow35A5SVbki\n/QIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D6").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAoLxrBD+gRvDI8Fs+kjHY\nfJPVB9O5Dm0NCy2XdumIk1PkmG2EzXqi3pmZLiF84D08WEvcZkbOQGdlbdqHDI1j\nzgI3fDTLvi5rlUYdn1s87a/7n7A6TIxLTTni23McAORYbRNAUpLNRojTNrnMelsY\nTtXgm3DWSuRG2w2+1IhKiQaA8jR19meFOYlL6UF9omPkOHmYkApt6e+aksjigfIq\nyJ8FE7BjnrZ0KFNTpplhf035HsVJV+8sih0rzQW30iY2H9NStinRAg2DaXNpL+7D\nUZL97JQOYdzNmfIRhhWvhmjDFcX7vqaKa5/0Z8ghj+oZHi9Ljm9kAfUFqy0jFPUB\niwIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D7").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAsRZwLQqUvz0/8xWVPuSj\nQfJUEx7F4oQ9RkIETNvDzSofRsU7vMPzJZM8lXcuAaGApRz6aAqe+P4SbRswG51j\nhOkfpFn0Pa3mcqQb4guA/GJ+OU1QiNSbbAzKZ2//3jHnHOcWCQZIeZwNAqIKcwqz\nur+wr3+qLuAT3VLBLFZfbIZqOGusDRjmIH0AM+Jd2lt9pATvHEjYukY5+Ksby6Rm\nqakGHFQ12TScZxnmSPLQNFXCR9eK0U7sLhSasUudSFvgj25l3jCc6SK5SDnaz3UA\n/UmaJq/deiqAOU+J7McD2qu27sOsjdnngqwfpE53HpbZ6V4i2wiHvsnS/pQBVVow\nDwIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D8").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAn0LfUPRGKkDomVnhtfgk\nKnWF5QQlDqMBFNuTFOfn2Ffy6dhSZsAlru/I52YE3lee5uXbm9ap64ndtQrvqKO8\nlkLobkkEKI/+OqCdRRPvUouMtPINyVGOujLvJb31AG87qW1Lb+hmTOxBES9dI9wB\nAEZluZsUKfL5Yt7HqUSrH50c8XRgobdzBqazw1T82AC9jFPd4Uu1pxLR3VThVMEc\nz/6yJhTB/m+EEnYEIsTrmXTMS+W5mmB/JRxv65+mKM6zT7QjTyRq3OtsP20Vkzw9\nHXmllzbaMbqK7pUlgRF0JkWdyp6XIkM8SaHcNcCrNI8VWSKOcJFwNZ3tGwLH7UMD\nRQIDAQAB\n-----END PUBLIC KEY-----\n'"
$ws.Range("D9").Value = "b'-----BEGIN PUBLIC KEY-----\nMIIBIjANBgkqhkiG9w0BAQEFAAOCAQ8AMIIBCgKCAQEAnApXEIsh3H4fcHGbtBf4\nB8lc6VF1Jk8K1S2tIZ+WGKfXEWpKiseHaMpj7d+h/3ETJR30dYQ3o3ti+wOCScJw\npVmYCi3+VbikJHZU8fwK9w35BsGTgoSFBBCb38I9hhAfOUOcekn2X7IegILaS4d7\n1DLFh6TTIkx7zAO/CVIQWH94fJB+tA0UCF7D2yT++/I69u9A6j9eP1jwpqedsFkw\ng7heiXl64x3aYKnt4HIuAhwQPKaCMA/k84jTGrhv4BqoiaoAooI6TsuQtwQmHMgn\nMzZPsR9JIuYX92wfGZ/Xb2rSP8cUU4TKyS+z8MK7wgNGGdYtYf/d7kHIHD/VwoJ6\nrwIDAQAB\n-----END PUBLIC KEY-----\n'"
